$d = $word.ActiveDocument

# The commit re-saved the document after list bookkeeping changes elsewhere
# in the file (new pH / soil organic content list items were added), which
# caused Word to mint fresh internal nsid GUIDs for several pre-existing
# multilevel list definitions in word/numbering.xml. Reproduce that exact
# set of nsid substitutions by round-tripping the package's flat-OPC XML
# through Document.WordOpenXML and rewriting just those attribute values.

$xml = $d.WordOpenXML

$map = @{
    "1b4d6395" = "eaeec97c"
    "dbd82f02" = "c56c3b59"
    "fef6db97" = "4d1542ac"
    "d974ebc8" = "73199394"
    "699fb8bb" = "8ecdbb64"
    "ab688dd0" = "93e23ec2"
    "986017b4" = "514a81da"
}

foreach ($old in $map.Keys) {
    $new = $map[$old]
    $xml = $xml.Replace('w:nsid w:val="' + $old + '"', 'w:nsid w:val="' + $new + '"')
}

$d.WordOpenXML = $xml

Write-Output "nsid values updated"
